# Applies the "Rename refactor" commit to Test-Spreadsheet.xlsx:
#   - TestRecord!A10   43269   -> 43270
#   - TestRecord!B10   128.34  -> 129.54
#   - Budget Out!C9    97.42   -> 98.62
#   - Expected Out!B9  1355.36 -> 1356.56
#   - Expected Out!B11 435.22  -> 436.42
#   - Expected Out!B1  =SUM(B2:B295) recalculates automatically from the
#     B9/B11 edits above (9680.195 -> 9682.595)
#   - Two shared strings each gain one extra trailing "z":
#       "Description007" + zzzz...  (50 z's, total len 64)  -> 51 z's (len 65)
#       "some test text" + zzzz...  (90 z's, total len 104) -> 91 z's (len 105)

$wb = $excel.ActiveWorkbook

$testRecord = $wb.Worksheets.Item("TestRecord")
$testRecord.Range("A10").Value = 43270
$testRecord.Range("B10").Value = 129.54

$budgetOut = $wb.Worksheets.Item("Budget Out")
$budgetOut.Range("C9").Value = 98.62

$expectedOut = $wb.Worksheets.Item("Expected Out")
$expectedOut.Range("B9").Value = 1356.56
$expectedOut.Range("B11").Value = 436.42

# Shared-string edits: append one extra "z" to each of the two long
# padded strings wherever they occur in the workbook. Build the "z" runs
# with an explicit loop (the "z" * n repeat-operator idiom isn't reliable
# here) so the lengths are exactly right.
# "Description007" + 50 z's (len 64) -> + 51 z's (len 65)
# "some test text" + 90 z's (len 104) -> + 91 z's (len 105)
$zRun50 = ""
for ($i = 0; $i -lt 50; $i++) { $zRun50 = $zRun50 + "z" }
$zRun51 = $zRun50 + "z"

$zRun90 = ""
for ($i = 0; $i -lt 90; $i++) { $zRun90 = $zRun90 + "z" }
$zRun91 = $zRun90 + "z"

$oldDescription = "Description007" + $zRun50
$newDescription = "Description007" + $zRun51
$oldSomeTestText = "some test text" + $zRun90
$newSomeTestText = "some test text" + $zRun91

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $current = $cell.Value2
        if ($current -eq $oldDescription) {
            $cell.Value = $newDescription
        } elseif ($current -eq $oldSomeTestText) {
            $cell.Value = $newSomeTestText
        }
    }
}
